$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.019.36'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '2.046.83'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'246.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").Value = "'57.50"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'0.0774"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Value = "'0.892"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +12.45%  '
$ws.Range("D14").Value = '2.345.82'
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("E15").Value = '  +1.67%  '
$ws.Range("D16").Value = '2.043.07'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = "'18.27"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +12.20%  '
$ws.Range("D18").Value = '36.956.21'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = "'74.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '0.0₃0899'
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Value = "'235.96"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  +4.11%  '
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("D26").Value = "'169.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("D27").Value = "'2.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("E29").Value = '  +15.37%  '
$ws.Range("D31").Value = "'1.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.06%  '
$ws.Range("D32").Value = "'4.77"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +7.73%  '
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").Value = "'0.0870"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.76%  '
$ws.Range("E36").Value = '  +2.39%  '
$ws.Range("E37").Value = '  +4.86%  '
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("D39").Value = "'3.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").Value = "'5.10"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.28%  '
$ws.Range("D41").Value = "'0.0991"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.78%  '
$ws.Range("D42").Value = "'0.0222"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").Value = "'97.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.10%  '
$ws.Range("D45").Value = "'17.01"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("D46").Value = '1.295.71'
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("E47").Value = '  -3.25%  '
$ws.Range("D49").Value = "'6.79"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("D50").Value = "'3.67"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.78%  '
$ws.Range("D51").Value = '2.230.08'
$ws.Range("E51").Value = '  -0.34%  '
